# Scheduled-runner data refresh: update Leve market-price / profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) across the eight crafting-job sheets with freshly
# pulled market values. Cells whose computed profit is now blank are
# cleared so the workbook matches Excel's "no value" representation
# (no cached zero left behind); cells that now have a value are created
# the same way a plain Range.Value assignment would.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 489.66666
$ws.Range("I12").Value = 584.5
$ws.Range("K12").Value = 584.5
$ws.Range("M12").Value = -414.5

$ws.Range("H19").Value = 638.13336
$ws.Range("I19").Value = 559.5
$ws.Range("J19").Value = 795.4
$ws.Range("K19").Value = 559.5
$ws.Range("L19").Value = 795.4
$ws.Range("M19").Value = -384.5
$ws.Range("N19").Value = -1145.4

$ws.Range("H33").Value = 155.07692
$ws.Range("I33").Value = 155.07692
$ws.Range("K33").Value = 155.07692
$ws.Range("M33").Value = 73.92308

$ws.Range("H80").Value = 222.81818
$ws.Range("J80").Value = 184.08333
$ws.Range("L80").Value = 552.24999
$ws.Range("N80").Value = -2548.24999

$ws.Range("H83").Value = 222.81818
$ws.Range("J83").Value = 184.08333
$ws.Range("L83").Value = 1656.74997
$ws.Range("N83").Value = -11640.74997

$ws.Range("H135").Value = 804.25714
$ws.Range("I135").Value = 603.2692
$ws.Range("K135").Value = 5429.422799999999
$ws.Range("M135").Value = -2894.422799999999

$ws.Range("H138").Value = 7352
$ws.Range("J138").Value = 6347.9473
$ws.Range("L138").Value = 19043.8419
$ws.Range("N138").Value = -29323.8419

$ws.Range("H141").Value = 3064.5386
$ws.Range("I141").Value = 2403
$ws.Range("J141").Value = 5269.6665
$ws.Range("K141").Value = 7209
$ws.Range("L141").Value = 15808.9995
$ws.Range("M141").Value = -2029
$ws.Range("N141").Value = -26168.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2666.3333
$ws.Range("I2").Value = 1499.5
$ws.Range("K2").Value = 1499.5
$ws.Range("M2").Value = -1386.5

$ws.Range("H32").Value = 4738.4443
$ws.Range("J32").Value = 10530
$ws.Range("L32").Value = 10530
$ws.Range("N32").Value = -11104

$ws.Range("H45").Value = 1925.25
$ws.Range("I45").Value = 1925.25
$ws.Range("K45").Value = 1925.25
$ws.Range("M45").Value = -1548.25

$ws.Range("H110").Value = 8354.200000000001
$ws.Range("I110").Value = 6134.7144
$ws.Range("K110").Value = 6134.7144
$ws.Range("M110").Value = -4089.7144

$ws.Range("H116").Value = 2666.3333
$ws.Range("I116").Value = 1499.5
$ws.Range("K116").Value = 1499.5
$ws.Range("M116").Value = 794.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2666.3333
$ws.Range("I3").Value = 1499.5
$ws.Range("K3").Value = 1499.5
$ws.Range("M3").Value = -1385.5

$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = $null

$ws.Range("H20").Value = 8785.571
$ws.Range("I20").Value = 9416.5
$ws.Range("K20").Value = 9416.5
$ws.Range("M20").Value = -9169.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4236.9805
$ws.Range("I31").Value = 2517.48
$ws.Range("J31").Value = 5890.346
$ws.Range("K31").Value = 2517.48
$ws.Range("L31").Value = 5890.346
$ws.Range("M31").Value = -2222.48
$ws.Range("N31").Value = -6480.346

$ws.Range("H34").Value = 4236.9805
$ws.Range("I34").Value = 2517.48
$ws.Range("J34").Value = 5890.346
$ws.Range("K34").Value = 2517.48
$ws.Range("L34").Value = 5890.346
$ws.Range("M34").Value = -2315.48
$ws.Range("N34").Value = -6294.346

$ws.Range("H58").Value = 3909.35
$ws.Range("I58").Value = 1648.2858
$ws.Range("K58").Value = 1648.2858
$ws.Range("M58").Value = -1445.2858

$ws.Range("H99").Value = 13386.826
$ws.Range("I99").Value = 9388.700000000001
$ws.Range("K99").Value = 9388.700000000001
$ws.Range("M99").Value = -7890.700000000001

$ws.Range("H126").Value = 13386.826
$ws.Range("I126").Value = 9388.700000000001
$ws.Range("K126").Value = 28166.1
$ws.Range("M126").Value = -25696.1

$ws.Range("H136").Value = 3909.35
$ws.Range("I136").Value = 1648.2858
$ws.Range("K136").Value = 4944.857400000001
$ws.Range("M136").Value = -2394.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 459.72726
$ws.Range("I12").Value = 439.66666
$ws.Range("J12").Value = 471.1905
$ws.Range("K12").Value = 1318.99998
$ws.Range("L12").Value = 1413.5715
$ws.Range("M12").Value = -1145.99998
$ws.Range("N12").Value = -1759.5715

$ws.Range("H22").Value = 5750
$ws.Range("J22").Value = 5750
$ws.Range("L22").Value = 17250
$ws.Range("N22").Value = -17588

$ws.Range("H23").Value = 167
$ws.Range("I23").Value = 199
$ws.Range("J23").Value = 135
$ws.Range("K23").Value = 597
$ws.Range("L23").Value = 405
$ws.Range("M23").Value = -362
$ws.Range("N23").Value = -875

$ws.Range("H27").Value = 5750
$ws.Range("J27").Value = 5750
$ws.Range("L27").Value = 17250
$ws.Range("N27").Value = -17454

$ws.Range("H122").Value = 124
$ws.Range("J122").Value = 48
$ws.Range("L122").Value = 432
$ws.Range("N122").Value = -5332

$ws.Range("H129").Value = 2194.5
$ws.Range("I129").Value = 2194.5
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 6583.5
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -1583.5
$ws.Range("N129").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 415.27777
$ws.Range("J2").Value = 633.36365
$ws.Range("L2").Value = 633.36365
$ws.Range("N2").Value = -859.36365

$ws.Range("H41").Value = 955.5
$ws.Range("I41").Value = 935
$ws.Range("J41").Value = 962.3333
$ws.Range("K41").Value = 935
$ws.Range("L41").Value = 962.3333
$ws.Range("M41").Value = -580
$ws.Range("N41").Value = -1672.3333

$ws.Range("H70").Value = 6617
$ws.Range("I70").Value = 6186.909
$ws.Range("K70").Value = 6186.909
$ws.Range("M70").Value = -5916.909

$ws.Range("H73").Value = 6617
$ws.Range("I73").Value = 6186.909
$ws.Range("K73").Value = 6186.909
$ws.Range("M73").Value = -5250.909

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null

$ws.Range("H126").Value = 4764.25
$ws.Range("J126").Value = 5169
$ws.Range("L126").Value = 15507
$ws.Range("N126").Value = -20447

$ws.Range("H132").Value = 2552.7083
$ws.Range("I132").Value = 1952.2
$ws.Range("K132").Value = 5856.6
$ws.Range("M132").Value = -3326.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 29291.5
$ws.Range("J20").Value = 28906
$ws.Range("L20").Value = 28906
$ws.Range("N20").Value = -29358

$ws.Range("H40").Value = 3197
$ws.Range("I40").Value = 3216.5454
$ws.Range("K40").Value = 3216.5454
$ws.Range("M40").Value = -3080.5454

$ws.Range("H55").Value = 1359.75
$ws.Range("I55").Value = 2119.5
$ws.Range("J55").Value = 600
$ws.Range("K55").Value = 2119.5
$ws.Range("L55").Value = 600
$ws.Range("M55").Value = -1946.5
$ws.Range("N55").Value = -946

$ws.Range("H82").Value = 2247.111
$ws.Range("I82").Value = 2518.5715
$ws.Range("J82").Value = 1297
$ws.Range("K82").Value = 2518.5715
$ws.Range("L82").Value = 1297
$ws.Range("M82").Value = -2157.5715
$ws.Range("N82").Value = -2019

$ws.Range("H85").Value = 2247.111
$ws.Range("I85").Value = 2518.5715
$ws.Range("J85").Value = 1297
$ws.Range("K85").Value = 2518.5715
$ws.Range("L85").Value = 1297
$ws.Range("M85").Value = -1270.5715
$ws.Range("N85").Value = -3793

$ws.Range("H136").Value = 5964
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 74498
$ws.Range("I60").Value = 110000
$ws.Range("K60").Value = 110000
$ws.Range("M60").Value = -109178

$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
